# Apply the two title-text edits described by the diff.
$p = $ppt.ActivePresentation

# -----------------------------------------------------------------------
# Slide 10 ("On next steps of meta-analysis" title):
#   "On next steps of meta-analysis" -> "On " / "next " / "steps" (3 runs)
# -----------------------------------------------------------------------
$s10  = $p.Slides.Item(10)
$sh10 = $s10.Shapes.Title
$tr10 = $sh10.TextFrame.TextRange

# Drop the trailing " of meta-analysis" (18 chars starting right after "steps").
$tail = $tr10.Characters(14, 18)
$tail.Delete()

# Text is now "On next steps"; split it into three separate runs by
# re-assigning the "next " and "steps" sub-ranges in place.
$run2 = $tr10.Characters(4, 5)
$run2.Text = "next "

$run3 = $tr10.Characters(9, 5)
$run3.Text = "steps"

# -----------------------------------------------------------------------
# Slide 8 ("cis-/trans- classification (INTERVAL/)" title):
#   "cis-/trans- " + "classification (INTERVAL/)" -> single merged run
# -----------------------------------------------------------------------
$s8  = $p.Slides.Item(8)
$sh8 = $s8.Shapes.Title
$tr8 = $sh8.TextFrame.TextRange

$firstRun  = $tr8.Characters(1, 12)
$secondRun = $tr8.Characters(13, 27)
$secondRun.Delete()
$firstRun.Text = "cis-/trans- classification (INTERVAL/)"
